# This patch removes internal rounding of floating point values used
# throughout the simulation. As a result, several previously rounded /
# integer reference values in this sheet must be replaced with their
# full-precision floating point equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = 17582.04134015827
    "C2"  = 37516.42916546755
    "D2"  = 73088.57832048151
    "E2"  = 110628.2963862276
    "B3"  = 180595.7446954236
    "C3"  = 335812.3065018094
    "D3"  = 407641.3939537266
    "E3"  = 453676.8843593199
    "B4"  = 19621.97193034408
    "C4"  = 37299.8589552691
    "D4"  = 59577.15749443948
    "E4"  = 77209.58902615773
    "B6"  = 105940.9511495386
    "C6"  = 135029.4366023016
    "D6"  = 124298.172958796
    "E6"  = 101988.9061478788
    "B7"  = 10963.85102564702
    "C7"  = 22019.45082437705
    "D7"  = 24319.82066029942
    "E7"  = 26523.24176013275
    "B9"  = 812445.7586883416
    "C9"  = 1305841.298497807
    "D9"  = 1771396.520995631
    "E9"  = 2185011.725311886
    "B12" = 900744.5176046697
    "C12" = 1064663.236149865
    "D12" = 948598.7011254284
    "E12" = 725999.6475283061
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
